$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Param1")

# Rename the "Animal  ID" header on Param1 to "Animal_ID"
$ws1.Range("B1").Value = "Animal_ID"

# Rename the Parameter1 time-point headers from dash to underscore
$ws1.Range("C1").Value = "Parameter1_0h"
$ws1.Range("D1").Value = "Parameter1_12h"

# Move the active selection on Param1 to B4 (was C13)
$ws1.Activate()
$ws1.Range("B4").Select()
